# Task 3: add a "Ratio with previous" column (D) to each of the four lab4
# worksheets. Column D holds, for rows 3-7, the ratio of the current row's
# runtime (column C) to the previous row's runtime: =C<r-1>/C<r>. Row 1 gets
# a bold header label; row 2 is left untouched (no "previous" row to ratio
# against).

$wb = $excel.ActiveWorkbook

$sheetNames = @("TwoSum", "TwoSumFast", "ThreeSum", "ThreeSumFast")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Header cell, styled to match the existing bold headers in A1:C1.
    $ws.Range("D1").Value = "Ratio with previous"
    $ws.Range("D1").Font.Bold = $true

    # Ratio formulas for rows 3-7 (row 2 has no "previous" row).
    $ws.Range("D3").Formula = "=C2/C3"
    $ws.Range("D4").Formula = "=C3/C4"
    $ws.Range("D5").Formula = "=C4/C5"
    $ws.Range("D6").Formula = "=C5/C6"
    $ws.Range("D7").Formula = "=C6/C7"
}

# Restore per-sheet selection state to match what was left behind after
# filling in the new column on each tab.
$wsTwoSum = $wb.Worksheets.Item("TwoSum")
$wsTwoSum.Activate() | Out-Null
$wsTwoSum.Range("D1:D7").Select() | Out-Null

$wsTwoSumFast = $wb.Worksheets.Item("TwoSumFast")
$wsTwoSumFast.Activate() | Out-Null
$wsTwoSumFast.Range("D3").Select() | Out-Null

$wsThreeSum = $wb.Worksheets.Item("ThreeSum")
$wsThreeSum.Activate() | Out-Null
$wsThreeSum.Range("D1:D7").Select() | Out-Null

# ThreeSumFast is the sheet that was active/visible when the workbook was
# saved (it was scrolled down to row 7 beforehand); scroll back to the top
# and select D1:D7 there, leaving it as the active tab.
$wsThreeSumFast = $wb.Worksheets.Item("ThreeSumFast")
$wsThreeSumFast.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsThreeSumFast.Range("D1:D7").Select() | Out-Null
